# Apply changes described by the commit: add new "Length movement" sheet
# with final length-of-movement measurement data, positioned right after
# "Blad1", and make it the active/selected sheet. Also update the
# selection on "Blad1" to E23 and deselect its tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Move the active selection on the original sheet to E23 (matches target
# sheetView selection) before we switch the active tab away from it.
$null = $ws1.Range("E23").Select()

# Insert the new worksheet right after "Blad1".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Length movement"

$cellData = @(
    @(1, 1, 30),
    @(1, 2, "avg mean"),
    @(1, 3, "avg std"),
    @(1, 4, "max"),
    @(2, 1, 0),
    @(2, 2, 81.6),
    @(2, 3, 0.17),
    @(2, 4, 81.78),
    @(3, 1, 1.25),
    @(3, 2, 80.09),
    @(3, 3, 0.095),
    @(3, 4, 80.18),
    @(4, 1, 1),
    @(4, 2, 76.41),
    @(4, 3, 0.22),
    @(4, 4, 76.63),
    @(5, 1, 0.75),
    @(5, 2, 76.65),
    @(5, 3, 0.47),
    @(5, 4, 77.12),
    @(6, 1, 0.5),
    @(6, 2, 76.86),
    @(6, 3, 1.21),
    @(6, 4, 78.08),
    @(7, 1, "Solar"),
    @(7, 2, 80.72),
    @(7, 3, 0.49),
    @(7, 4, 81.16),
    @(9, 1, 50),
    @(9, 2, "avg mean"),
    @(9, 3, "avg std"),
    @(9, 4, "max"),
    @(10, 1, 0),
    @(10, 2, 79.62),
    @(10, 3, 0.47),
    @(10, 4, 80.09),
    @(11, 1, 1.25),
    @(11, 2, 79.91),
    @(11, 3, 0.34),
    @(11, 4, 80.26),
    @(12, 1, 1),
    @(12, 2, 80.35),
    @(12, 3, 0.17),
    @(13, 1, 0.75),
    @(13, 2, 77.95),
    @(13, 3, 0.153),
    @(13, 4, 78.1),
    @(14, 1, 0.5),
    @(14, 2, 79.2),
    @(14, 3, 0.22),
    @(14, 4, 79.2),
    @(15, 1, "Solar"),
    @(15, 2, 80.22),
    @(15, 3, 1.05),
    @(15, 4, 81.7),
    @(17, 1, 70),
    @(17, 2, "avg mean"),
    @(17, 3, "avg std"),
    @(17, 4, "max"),
    @(18, 1, 0),
    @(18, 2, 73.72),
    @(18, 3, 0),
    @(18, 4, 73.72),
    @(19, 1, 1.25),
    @(19, 2, 73.44),
    @(19, 3, 0.4),
    @(19, 4, 73.82),
    @(20, 1, 1),
    @(20, 2, 73.21),
    @(20, 3, 0.15),
    @(20, 4, 73.37),
    @(21, 1, 0.75),
    @(21, 2, 74.04),
    @(21, 3, 0.57),
    @(21, 4, 74.61),
    @(22, 1, 0.5),
    @(22, 2, 71.3),
    @(22, 3, 0.02),
    @(22, 4, 71.35),
    @(23, 1, "Solar"),
    @(23, 2, 81.4),
    @(23, 3, 0.38),
    @(23, 4, 81.92)
)

foreach ($item in $cellData) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws2.Cells.Item($r, $c).Value = $v
}

# Set the selection on the new sheet and make it the active sheet/tab,
# matching the target sheetView (activeCell E22, tabSelected).
$null = $ws2.Range("E22").Select()
$null = $ws2.Activate()
